$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Ana's age (B2) from 21 to 19
$ws.Range("B2").Value = 19

# Update Luis's province (C3) from "Alajuela" to "Cartago"
$ws.Range("C3").Value = "Cartago"

# Update the selection to C3 (as in the diff)
$ws.Range("C3").Select()
